$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 101: Edge of the Arcane | Cunning Craftsman's Tea
$ws.Range("H101").Value = 798
$ws.Range("I101").Value = 839.375
$ws.Range("K101").Value = 2518.125
$ws.Range("M101").Value = -896.125

# Row 125: Body over Mind | Grade 5 Dexterity Alkahest
$ws.Range("H125").Value = 1036
$ws.Range("I125").Value = 1239.3334
$ws.Range("K125").Value = 11154.0006
$ws.Range("M125").Value = -8694.000599999999

# Row 131: Mindful Study | Grade 5 Tincture of Mind
$ws.Range("H131").Value = 1131.5
$ws.Range("I131").Value = 1131.5
$ws.Range("K131").Value = 3394.5
$ws.Range("M131").Value = 1645.5

# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 30306706
$ws.Range("I132").Value = 38465784
$ws.Range("J132").Value = 1555.8572
$ws.Range("K132").Value = 115397352
$ws.Range("L132").Value = 4667.571599999999
$ws.Range("M132").Value = -115394822
$ws.Range("N132").Value = -9727.571599999999

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 3905.3438
$ws.Range("J138").Value = 3873.0833
$ws.Range("L138").Value = 11619.2499
$ws.Range("N138").Value = -21899.2499


$ws = $wb.Worksheets.Item("ARM")
# Row 5: The Alloyed Truth | Bronze Rivets
$ws.Range("H5").Value = 668.375
$ws.Range("I5").Value = 475
$ws.Range("K5").Value = 475
$ws.Range("M5").Value = -363

# Row 63: Rivets Run through It | Mythrite Rivets
$ws.Range("H63").Value = 1899.6
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

# Row 66: A Riveting Revival (L) | Mythrite Rivets
$ws.Range("H66").Value = 1899.6
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

# Row 88: The Mast Chance | Adamantite Rivets
$ws.Range("H88").Value = 784.5
$ws.Range("I88").Value = 474.5
$ws.Range("J88").Value = 862
$ws.Range("K88").Value = 474.5
$ws.Range("L88").Value = 862
$ws.Range("M88").Value = -68.5
$ws.Range("N88").Value = -1674

# Row 91: The Rose and the Riveter (L) | Adamantite Rivets
$ws.Range("H91").Value = 784.5
$ws.Range("I91").Value = 474.5
$ws.Range("J91").Value = 862
$ws.Range("K91").Value = 474.5
$ws.Range("L91").Value = 862
$ws.Range("M91").Value = 929.5
$ws.Range("N91").Value = -3670

# Row 97: Ore for Me | High Steel Ingot
$ws.Range("H97").Value = 976.64703
$ws.Range("I97").Value = 881.5
$ws.Range("K97").Value = 881.5
$ws.Range("M97").Value = -385.5

# Row 102: Smells of Rich Tama-hagane | Tama-hagane Ingot
$ws.Range("H102").Value = 5573.7856
$ws.Range("I102").Value = 4324.15
$ws.Range("J102").Value = 8697.875
$ws.Range("K102").Value = 4324.15
$ws.Range("L102").Value = 8697.875
$ws.Range("M102").Value = -2702.15
$ws.Range("N102").Value = -11941.875

# Row 124: Ace of Gloves | High Durium Gauntlets of Fending
$ws.Range("H124").Value = 23500
$ws.Range("J124").Value = 23500
$ws.Range("L124").Value = 23500
$ws.Range("N124").Value = -33320

# Row 130: A Gift of Gloves | Chondrite Gloves of Casting
$ws.Range("H130").Value = 49078
$ws.Range("J130").Value = 49078
$ws.Range("L130").Value = 49078
$ws.Range("N130").Value = -59118

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 3684.8445
$ws.Range("I132").Value = 3673.6099
$ws.Range("K132").Value = 11020.8297
$ws.Range("M132").Value = -8490.8297


$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences | Bronze Rivets
$ws.Range("H4").Value = 668.375
$ws.Range("I4").Value = 475
$ws.Range("K4").Value = 475
$ws.Range("M4").Value = -360

# Row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value = 12375234
$ws.Range("I86").Value = 27473.3
$ws.Range("K86").Value = 27473.3
$ws.Range("M86").Value = -26350.3

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value = 12375234
$ws.Range("I89").Value = 27473.3
$ws.Range("K89").Value = 137366.5
$ws.Range("M89").Value = -131750.5

# Row 99: Meddle in Metal | Oroshigane Ingot
$ws.Range("H99").Value = 1525
$ws.Range("I99").Value = 1525
$ws.Range("K99").Value = 1525
$ws.Range("M99").Value = -27


$ws = $wb.Worksheets.Item("CRP")
# Row 19: Shielding Sales | Square Ash Shield
$ws.Range("H19").Value = 4379.9
$ws.Range("I19").Value = 421.44446
$ws.Range("K19").Value = 421.44446
$ws.Range("M19").Value = -251.44446

# Row 22: Driving Up the Wall | Elm Lumber
$ws.Range("H22").Value = 322.91666
$ws.Range("J22").Value = 414.5
$ws.Range("L22").Value = 414.5
$ws.Range("N22").Value = -1114.5

# Row 24: What You Need | Square Ash Shield
$ws.Range("H24").Value = 4379.9
$ws.Range("I24").Value = 421.44446
$ws.Range("K24").Value = 421.44446
$ws.Range("M24").Value = -251.44446

# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 2213.1482
$ws.Range("I31").Value = 2028.9445
$ws.Range("K31").Value = 2028.9445
$ws.Range("M31").Value = -1733.9445

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 2213.1482
$ws.Range("I34").Value = 2028.9445
$ws.Range("K34").Value = 2028.9445
$ws.Range("M34").Value = -1826.9445

# Row 122: Timber of Tenkonto | Horse Chestnut Lumber
$ws.Range("H122").Value = 4087.7693
$ws.Range("I122").Value = 4592.9
$ws.Range("J122").Value = 2404
$ws.Range("K122").Value = 13778.7
$ws.Range("L122").Value = 7212
$ws.Range("M122").Value = -11328.7
$ws.Range("N122").Value = -12112


$ws = $wb.Worksheets.Item("CUL")
# Row 33: Cooking with Gas | Chicken Stock
$ws.Range("H33").Value = 132.16667
$ws.Range("I33").Value = 139.33333
$ws.Range("J33").Value = 125
$ws.Range("K33").Value = 835.9999799999999
$ws.Range("L33").Value = 750
$ws.Range("M33").Value = -552.9999799999999
$ws.Range("N33").Value = -1316

# Row 34: Fever Pitch | Chamomile Tea
$ws.Range("H34").Value = 63983
$ws.Range("J34").Value = 90538.164
$ws.Range("L34").Value = 271614.492
$ws.Range("N34").Value = -271782.492

# Row 39: Bloody Good Tart, This | Blood Currant Tart
$ws.Range("H39").Value = 7900
$ws.Range("J39").Value = 7750
$ws.Range("L39").Value = 23250
$ws.Range("N39").Value = -23838

# Row 47: Winter of Our Discontent | Mugwort Carp
$ws.Range("H47").Value = 4613.3335
$ws.Range("I47").Value = 4613.3335
$ws.Range("K47").Value = 13840.0005
$ws.Range("M47").Value = -13409.0005

# Row 55: Pagan Pastries | Pastry Fish
$ws.Range("H55").Value = 338.75
$ws.Range("I55").Value = 338.75
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 1016.25
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -839.25
$ws.Range("N55").ClearContents()

# Row 86: Let's Not Get Sappy | Birch Syrup
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# Row 89: Luxury Spillover (L) | Birch Syrup
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()


$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers | Copper Ingot
$ws.Range("H2").Value = 192.45833
$ws.Range("I2").Value = 177.29411
$ws.Range("K2").Value = 177.29411
$ws.Range("M2").Value = -64.29410999999999

# Row 80: Needs More Prayerbell | Hardsilver Ingot
$ws.Range("H80").Value = 3802.6
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 4253.25
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 4253.25
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = -6249.25

# Row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
$ws.Range("H83").Value = 3802.6
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 4253.25
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 21266.25
$ws.Range("M83").Value = -5008
$ws.Range("N83").Value = -31250.25

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 23811900
$ws.Range("I132").Value = 2751.1667
$ws.Range("J132").Value = 41668760
$ws.Range("K132").Value = 8253.500100000001
$ws.Range("L132").Value = 125006280
$ws.Range("M132").Value = -5723.500100000001
$ws.Range("N132").Value = -125011340


$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban | Leather
$ws.Range("H7").Value = 11511.8
$ws.Range("I7").Value = 11511.8
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 11511.8
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -11399.8
$ws.Range("N7").ClearContents()

# Row 126: Battered Books | Saiga Leather
$ws.Range("H126").Value = 11511.8
$ws.Range("I126").Value = 11511.8
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 34535.39999999999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -32065.39999999999
$ws.Range("N126").ClearContents()

# Row 128: Grips of Fear | Kumbhiraskin Gloves of the Black Griffin
$ws.Range("H128").Value = 78993.5
$ws.Range("J128").Value = 78993.5
$ws.Range("L128").Value = 78993.5
$ws.Range("N128").Value = -88953.5

# Row 130: Generous Soles | Ophiotauroskin Boots of Healing
$ws.Range("H130").Value = 75943
$ws.Range("J130").Value = 75943
$ws.Range("L130").Value = 75943
$ws.Range("N130").Value = -85983


$ws = $wb.Worksheets.Item("WVR")
# Row 107: Flax Wax | Bright Linen Yarn
$ws.Range("H107").Value = 304.08334
$ws.Range("I107").Value = 304.08334
$ws.Range("K107").Value = 912.2500200000001
$ws.Range("M107").Value = 1007.74998

# Row 113: A Tender Table | Pixie Floss
$ws.Range("H113").Value = 794
$ws.Range("I113").Value = 795.8570999999999
$ws.Range("K113").Value = 2387.5713
$ws.Range("M113").Value = -217.5712999999996

# Row 122: Heavy Armoire | Dark Hempen Cloth
$ws.Range("H122").Value = 2849.5
$ws.Range("I122").Value = 2706.1333
$ws.Range("K122").Value = 8118.3999
$ws.Range("M122").Value = -5668.3999

# Row 126: A Polished Purchase | Snow Linen
$ws.Range("H126").Value = 3207258.2
$ws.Range("I126").Value = 5954438.5
$ws.Range("J126").Value = 2214.5
$ws.Range("K126").Value = 17863315.5
$ws.Range("L126").Value = 6643.5
$ws.Range("M126").Value = -17860845.5
$ws.Range("N126").Value = -11583.5

